$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("March_CourseList")

# Insert a new row at position 10, pushing the existing "CIA (morning)" row down to row 11
$ws.Rows(10).Insert()

# The pushed-down row (now row 11) keeps the old text; rename it to "CIA(Evening)"
$ws.Range("A11").Value = "CIA(Evening)"

# The newly inserted row 10 gets the restricted/shortened course name
$ws.Range("A10").Value = "CIA (Morning)"

# Update selection to match the new active cell
$ws.Range("A11").Select()
